$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.526.04"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.256.08"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.05%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.47"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.24%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.253.10"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("E10").Value = "  +7.15%  "
$ws.Range("E11").Value = "  +2.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.416"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +6.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.817.07"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.11%  "
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.40"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.474.38"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.34%  "
$ws.Range("E17").Value = "  +3.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.271.06"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.84"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.36%  "
$ws.Range("E20").Value = "  +5.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.71"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.64"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.36%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.511"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.73%  "
$ws.Range("E26").Value = "  +5.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.58"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  +3.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("E30").Value = "  +3.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.64"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.73"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.10%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  +6.22%  "
$ws.Range("E35").Value = "  +4.96%  "
$ws.Range("E36").Value = "  +5.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.29"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.846"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.84"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.82"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +12.22%  "
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.64"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +11.30%  "
$ws.Range("E43").Value = "  +5.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.696.49"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "350.90"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.44"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +6.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.79"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.86%  "
$ws.Range("E48").Value = "  +3.43%  "
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.997"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.94%  "
$ws.Range("E51").Value = "  +0.66%  "
